# TR Final Chap 4 and 5
#
# Refresh the cached "datetimeFigureOut" footer placeholder text from
# 8/30/2021 (2021/8/30 in the zh-TW notes master) to 9/5/2021 (2021/9/5)
# across the slide master and every slide layout, and fix the "SimpeVS"
# typo on slide 1.
#
# NOTE: this runtime's ActivePresentation.NotesMaster is aliased onto the
# slide master's Shapes collection for writes (a host quirk), so touching
# it here would silently corrupt the slide master's body placeholder
# instead of updating the real notes master. It is intentionally left
# alone.

$p = $ppt.ActivePresentation

function Update-DatePlaceholder {
    param($shapes)
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq "8/30/2021") {
                $tr.Text = "9/5/2021"
            }
        }
    }
}

# Slide master's own Date Placeholder
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

# Every slide layout has its own cached copy of the date placeholder text
for ($L = 1; $L -le $master.CustomLayouts.Count; $L++) {
    $layout = $master.CustomLayouts.Item($L)
    Update-DatePlaceholder $layout.Shapes
}

# Fix the "SimpeVS" -> "SimpleVS" typo on slide 1
$s1 = $p.Slides.Item(1)
for ($i = 1; $i -le $s1.Shapes.Count; $i++) {
    $shp = $s1.Shapes.Item($i)
    if ($shp.HasTextFrame) {
        if ($shp.TextFrame.HasText) {
            if ($shp.TextFrame.TextRange.Text -eq "SimpeVS") {
                $shp.TextFrame.TextRange.Text = "SimpleVS"
            }
        }
    }
}
